# Generate Report for Handoff
# A new handoff was recorded for file "8f1ee2d7-3d3c-4371-b975-19413b8bf9ee.md"
# (which "ba234996-fc6a-4969-8532-6f2097065cec.md" depends on), so the
# "Latest Handoff Datetime" column (D) is updated for both of those rows
# (row 12 = 8f1ee2d7..., row 14 = ba234996...) on each localized-language
# status sheet.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D12").Value = "2016-03-09 10:29:02"
$zhcn.Range("D14").Value = "2016-03-09 10:29:02"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D12").Value = "2016-03-09 10:29:06"
$dede.Range("D14").Value = "2016-03-09 10:29:06"
